$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark left over in the "Fjhty" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append a new run containing "kyfiufiutf" right after the run that holds "2" ---
# Find the paragraph whose entire text is "2" (it is unique in this document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $pr = $p.Range
    if ($pr.Text.TrimEnd() -eq "2") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pr = $target.Range
    # Exclude the trailing paragraph mark from the range.
    $pr.End = $pr.End - 1

    # Collapse to the very end of the existing text and insert the new text there.
    $ins = $pr.Duplicate
    $ins.Collapse(0)
    $ins.InsertAfter("kyfiufiutf")

    # The engine merges adjacent runs that end up with byte-identical rPr, so the
    # freshly inserted text would otherwise be folded back into the "2" run. Force
    # the new text to become (and remain) its own <w:r> by toggling a run-level
    # formatting flag on/off around it -- this mirrors the two separate <w:r>
    # elements (both carrying only <w:lang w:val="en-US"/>) seen in the target
    # document while still leaving the run's effective formatting unchanged.
    $newRange = $d.Range($ins.Start, $ins.Start + 10)
    $newRange.Bold = 1
    $newRange.Bold = 0
}
